$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'300.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.88%"
$ws.Range("E2").Style = "Normal"
$ws.Range("E3").Value = "'-1.80%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.111"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-3.07%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07359"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-1.80%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.348"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'57.30%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.951"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'1.23%"
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'-0.79%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9182"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.15%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1709"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'1.13%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07588"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.62%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.08103"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.04%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03018"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.62%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.09925"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.28%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001495"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.07%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.006090"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-5.72%"
$ws.Range("E16").Style = "Normal"
$ws.Range("E17").Value = "'0.10%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'-0.27%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'-0.56%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'-0.52%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.652"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'3.16%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.04632"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.78%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.1567"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-3.32%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'0.87%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004478"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.83%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.0001299"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.05%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'-3.26%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01720"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-0.02%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04521"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.91%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007196"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'3.86%"
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "'-0.15%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002238"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'1.46%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.01073"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-15.87%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006271"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'1.62%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.009992"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-33.33%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.8085"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-56.66%"
$ws.Range("E47").Style = "Normal"
